$p = $ppt.ActivePresentation

# --- Slide 34: merge "Magnetometer & accelerometer..." runs into one ---
$s34 = $p.Slides.Item(34)
$shp34 = $s34.Shapes.Item(2)
$shp34.TextFrame.TextRange.Text = "Doesn’t incorporate gyro dataAssumes that the acceleration vector always points downDoesn’t fix the data rate as a constantMagnetometer & accelerometer need calibration for best performance."

# --- Slide 37: add "Arduino Cookbook (Margolis)" paragraph ---
$s37 = $p.Slides.Item(37)
$shp37 = $s37.Shapes.Item(3)
$tr37 = $shp37.TextFrame.TextRange
$lastNonEmptyPara = $tr37.Paragraphs(4, 1)
$inserted = $lastNonEmptyPara.InsertAfter("`rArduino Cookbook ")
$newPara = $tr37.Paragraphs(5, 1)
$inserted2 = $newPara.InsertAfter("(Margolis)")
$trailingEmptyPara = $tr37.Paragraphs(6, 1)
$trailingEmptyPara.Delete()
